$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value2 = 4708.6523
$ws.Cells.Item(28, 9).Value2 = 6229.4116
$ws.Cells.Item(28, 10).Value2 = 399.83334
$ws.Cells.Item(28, 11).Value2 = 6229.4116
$ws.Cells.Item(28, 12).Value2 = 399.83334
$ws.Cells.Item(28, 13).Value2 = -5744.4116
$ws.Cells.Item(28, 14).Value2 = -1369.83334
# Row 38
$ws.Cells.Item(38, 8).Value2 = 276
$ws.Cells.Item(38, 9).Value2 = 117.09091
$ws.Cells.Item(38, 10).Value2 = 1150
$ws.Cells.Item(38, 11).Value2 = 351.27273
$ws.Cells.Item(38, 12).Value2 = 3450
$ws.Cells.Item(38, 13).Value2 = 20.72727000000003
$ws.Cells.Item(38, 14).Value2 = -4194
# Row 62
$ws.Cells.Item(62, 8).Value2 = 3430.4443
$ws.Cells.Item(62, 9).Value2 = 3244.25
$ws.Cells.Item(62, 10).Value2 = 3579.4
$ws.Cells.Item(62, 11).Value2 = 3244.25
$ws.Cells.Item(62, 12).Value2 = 3579.4
$ws.Cells.Item(62, 13).Value2 = -2620.25
$ws.Cells.Item(62, 14).Value2 = -4827.4
# Row 64
$ws.Cells.Item(64, 8).Value2 = 3442.027
$ws.Cells.Item(64, 9).Value2 = 3529.3845
$ws.Cells.Item(64, 10).Value2 = 3394.7083
$ws.Cells.Item(64, 11).Value2 = 3529.3845
$ws.Cells.Item(64, 12).Value2 = 3394.7083
$ws.Cells.Item(64, 13).Value2 = -3281.3845
$ws.Cells.Item(64, 14).Value2 = -3890.7083
# Row 65
$ws.Cells.Item(65, 8).Value2 = 3430.4443
$ws.Cells.Item(65, 9).Value2 = 3244.25
$ws.Cells.Item(65, 10).Value2 = 3579.4
$ws.Cells.Item(65, 11).Value2 = 16221.25
$ws.Cells.Item(65, 12).Value2 = 17897
$ws.Cells.Item(65, 13).Value2 = -13101.25
$ws.Cells.Item(65, 14).Value2 = -24137
# Row 67
$ws.Cells.Item(67, 8).Value2 = 3442.027
$ws.Cells.Item(67, 9).Value2 = 3529.3845
$ws.Cells.Item(67, 10).Value2 = 3394.7083
$ws.Cells.Item(67, 11).Value2 = 3529.3845
$ws.Cells.Item(67, 12).Value2 = 3394.7083
$ws.Cells.Item(67, 13).Value2 = -2671.3845
$ws.Cells.Item(67, 14).Value2 = -5110.7083
# Row 74
$ws.Cells.Item(74, 8).Value2 = 3038.75
$ws.Cells.Item(74, 9).Value2 = 3038.75
$ws.Cells.Item(74, 11).Value2 = 3038.75
$ws.Cells.Item(74, 13).Value2 = -2102.75
# Row 77
$ws.Cells.Item(77, 8).Value2 = 3038.75
$ws.Cells.Item(77, 9).Value2 = 3038.75
$ws.Cells.Item(77, 11).Value2 = 15193.75
$ws.Cells.Item(77, 13).Value2 = -10513.75
# Row 100
$ws.Cells.Item(100, 8).Value2 = 2220.625
$ws.Cells.Item(100, 9).Value2 = 2402.1428
$ws.Cells.Item(100, 10).Value2 = 950
$ws.Cells.Item(100, 11).Value2 = 2402.1428
$ws.Cells.Item(100, 12).Value2 = 950
$ws.Cells.Item(100, 13).Value2 = -1861.1428
$ws.Cells.Item(100, 14).Value2 = -2032
# Row 103
$ws.Cells.Item(103, 8).Value2 = 1510
$ws.Cells.Item(103, 10).Value2 = 1787.5
$ws.Cells.Item(103, 12).Value2 = 5362.5
$ws.Cells.Item(103, 14).Value2 = -6534.5
# Row 111
$ws.Cells.Item(111, 8).Value2 = 4025.8125
$ws.Cells.Item(111, 10).Value2 = 5029.8
$ws.Cells.Item(111, 12).Value2 = 15089.4
$ws.Cells.Item(111, 14).Value2 = -21223.4
# Row 112
$ws.Cells.Item(112, 8).Value2 = 2337.5588
$ws.Cells.Item(112, 10).Value2 = 2410.25
$ws.Cells.Item(112, 12).Value2 = 7230.75
$ws.Cells.Item(112, 14).Value2 = -9446.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value2 = 3279.9
$ws.Cells.Item(32, 9).Value2 = 3766.3513
$ws.Cells.Item(32, 11).Value2 = 3766.3513
$ws.Cells.Item(32, 13).Value2 = -3479.3513
# Row 45
$ws.Cells.Item(45, 8).Value2 = 1988.4
$ws.Cells.Item(45, 9).Value2 = 1855.2
$ws.Cells.Item(45, 11).Value2 = 1855.2
$ws.Cells.Item(45, 13).Value2 = -1478.2
# Row 61
$ws.Cells.Item(61, 8).Value2 = 932.0417
$ws.Cells.Item(61, 9).Value2 = 958.913
$ws.Cells.Item(61, 11).Value2 = 958.913
$ws.Cells.Item(61, 13).Value2 = -746.913
# Row 63
$ws.Cells.Item(63, 8).Value2 = 166668670
$ws.Cells.Item(63, 9).Value2 = 2500
$ws.Cells.Item(63, 10).Value2 = 500001000
$ws.Cells.Item(63, 11).Value2 = 2500
$ws.Cells.Item(63, 12).Value2 = 500001000
$ws.Cells.Item(63, 13).Value2 = -1814
$ws.Cells.Item(63, 14).Value2 = -500002372
# Row 66
$ws.Cells.Item(66, 8).Value2 = 166668670
$ws.Cells.Item(66, 9).Value2 = 2500
$ws.Cells.Item(66, 10).Value2 = 500001000
$ws.Cells.Item(66, 11).Value2 = 12500
$ws.Cells.Item(66, 12).Value2 = 2500005000
$ws.Cells.Item(66, 13).Value2 = -9068
$ws.Cells.Item(66, 14).Value2 = -2500011864
# Row 122
$ws.Cells.Item(122, 8).Value2 = 1422.8462
$ws.Cells.Item(122, 9).Value2 = 1374.75
$ws.Cells.Item(122, 10).Value2 = 2000
$ws.Cells.Item(122, 11).Value2 = 4124.25
$ws.Cells.Item(122, 12).Value2 = 6000
$ws.Cells.Item(122, 13).Value2 = -1674.25
$ws.Cells.Item(122, 14).Value2 = -10900
# Row 132
$ws.Cells.Item(132, 8).Value2 = 2344.0977
$ws.Cells.Item(132, 9).Value2 = 2086
$ws.Cells.Item(132, 11).Value2 = 6258
$ws.Cells.Item(132, 13).Value2 = -3728
# Row 136
$ws.Cells.Item(136, 8).Value2 = 932.0417
$ws.Cells.Item(136, 9).Value2 = 958.913
$ws.Cells.Item(136, 11).Value2 = 2876.739
$ws.Cells.Item(136, 13).Value2 = -326.739

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value2 = 930.56525
$ws.Cells.Item(31, 9).Value2 = 727.6957
$ws.Cells.Item(31, 11).Value2 = 727.6957
$ws.Cells.Item(31, 13).Value2 = -432.6957
# Row 34
$ws.Cells.Item(34, 8).Value2 = 930.56525
$ws.Cells.Item(34, 9).Value2 = 727.6957
$ws.Cells.Item(34, 11).Value2 = 727.6957
$ws.Cells.Item(34, 13).Value2 = -525.6957
# Row 132
$ws.Cells.Item(132, 8).Value2 = 7697
$ws.Cells.Item(132, 9).Value2 = 9623.214
$ws.Cells.Item(132, 11).Value2 = 28869.642
$ws.Cells.Item(132, 13).Value2 = -26339.642

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Cells.Item(13, 8).Value2 = 599.25
$ws.Cells.Item(13, 9).Value2 = 465.66666
$ws.Cells.Item(13, 11).Value2 = 1396.99998
$ws.Cells.Item(13, 13).Value2 = -1228.99998
# Row 64
$ws.Cells.Item(64, 8).Value2 = 4171
$ws.Cells.Item(64, 9).Value2 = 1066.6666
$ws.Cells.Item(64, 10).Value2 = 4688.3887
$ws.Cells.Item(64, 11).Value2 = 3199.9998
$ws.Cells.Item(64, 12).Value2 = 14065.1661
$ws.Cells.Item(64, 13).Value2 = -2929.9998
$ws.Cells.Item(64, 14).Value2 = -14605.1661
# Row 67
$ws.Cells.Item(67, 8).Value2 = 4171
$ws.Cells.Item(67, 9).Value2 = 1066.6666
$ws.Cells.Item(67, 10).Value2 = 4688.3887
$ws.Cells.Item(67, 11).Value2 = 3199.9998
$ws.Cells.Item(67, 12).Value2 = 14065.1661
$ws.Cells.Item(67, 13).Value2 = -2263.9998
$ws.Cells.Item(67, 14).Value2 = -15937.1661
# Row 68
$ws.Cells.Item(68, 8).Value2 = 1746.7106
$ws.Cells.Item(68, 10).Value2 = 2029.1
$ws.Cells.Item(68, 12).Value2 = 6087.299999999999
$ws.Cells.Item(68, 14).Value2 = -7709.299999999999
# Row 70
$ws.Cells.Item(70, 8).Value2 = 5340
$ws.Cells.Item(70, 9).Value2 = 2000
$ws.Cells.Item(70, 10).Value2 = 5711.1113
$ws.Cells.Item(70, 11).Value2 = 6000
$ws.Cells.Item(70, 12).Value2 = 17133.3339
$ws.Cells.Item(70, 13).Value2 = -5685
$ws.Cells.Item(70, 14).Value2 = -17763.3339
# Row 71
$ws.Cells.Item(71, 8).Value2 = 1746.7106
$ws.Cells.Item(71, 10).Value2 = 2029.1
$ws.Cells.Item(71, 12).Value2 = 18261.9
$ws.Cells.Item(71, 14).Value2 = -26373.9
# Row 73
$ws.Cells.Item(73, 8).Value2 = 5340
$ws.Cells.Item(73, 9).Value2 = 2000
$ws.Cells.Item(73, 10).Value2 = 5711.1113
$ws.Cells.Item(73, 11).Value2 = 6000
$ws.Cells.Item(73, 12).Value2 = 17133.3339
$ws.Cells.Item(73, 13).Value2 = -4908
$ws.Cells.Item(73, 14).Value2 = -19317.3339
# Row 109
$ws.Cells.Item(109, 8).Value2 = 127669.25
$ws.Cells.Item(109, 9).Value2 = 334118
$ws.Cells.Item(109, 11).Value2 = 1002354
$ws.Cells.Item(109, 13).Value2 = -1001314

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value2 = 18754322
$ws.Cells.Item(70, 9).Value2 = 19234594
$ws.Cells.Item(70, 10).Value2 = 18186730
$ws.Cells.Item(70, 11).Value2 = 19234594
$ws.Cells.Item(70, 12).Value2 = 18186730
$ws.Cells.Item(70, 13).Value2 = -19234324
$ws.Cells.Item(70, 14).Value2 = -18187270
# Row 73
$ws.Cells.Item(73, 8).Value2 = 18754322
$ws.Cells.Item(73, 9).Value2 = 19234594
$ws.Cells.Item(73, 10).Value2 = 18186730
$ws.Cells.Item(73, 11).Value2 = 19234594
$ws.Cells.Item(73, 12).Value2 = 18186730
$ws.Cells.Item(73, 13).Value2 = -19233658
$ws.Cells.Item(73, 14).Value2 = -18188602
# Row 80
$ws.Cells.Item(80, 8).Value2 = 6044.3335
$ws.Cells.Item(80, 9).Value2 = 5999.5
$ws.Cells.Item(80, 11).Value2 = 5999.5
$ws.Cells.Item(80, 13).Value2 = -5001.5
# Row 83
$ws.Cells.Item(83, 8).Value2 = 6044.3335
$ws.Cells.Item(83, 9).Value2 = 5999.5
$ws.Cells.Item(83, 11).Value2 = 29997.5
$ws.Cells.Item(83, 13).Value2 = -25005.5
# Row 97
$ws.Cells.Item(97, 8).Value2 = 1196.8695
$ws.Cells.Item(97, 9).Value2 = 1587.6428
$ws.Cells.Item(97, 10).Value2 = 589
$ws.Cells.Item(97, 11).Value2 = 1587.6428
$ws.Cells.Item(97, 12).Value2 = 589
$ws.Cells.Item(97, 13).Value2 = -1091.6428
$ws.Cells.Item(97, 14).Value2 = -1581
# Row 132
$ws.Cells.Item(132, 8).Value2 = 2883.1428
$ws.Cells.Item(132, 9).Value2 = 2741.077
$ws.Cells.Item(132, 10).Value2 = 3114
$ws.Cells.Item(132, 11).Value2 = 8223.231
$ws.Cells.Item(132, 12).Value2 = 9342
$ws.Cells.Item(132, 13).Value2 = -5693.231
$ws.Cells.Item(132, 14).Value2 = -14402

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value2 = 1210.8889
$ws.Cells.Item(22, 9).Value2 = 619.8
$ws.Cells.Item(22, 10).Value2 = 1949.75
$ws.Cells.Item(22, 11).Value2 = 619.8
$ws.Cells.Item(22, 12).Value2 = 1949.75
$ws.Cells.Item(22, 13).Value2 = -324.8
$ws.Cells.Item(22, 14).Value2 = -2539.75
# Row 27
$ws.Cells.Item(27, 8).Value2 = 1210.8889
$ws.Cells.Item(27, 9).Value2 = 619.8
$ws.Cells.Item(27, 10).Value2 = 1949.75
$ws.Cells.Item(27, 11).Value2 = 619.8
$ws.Cells.Item(27, 12).Value2 = 1949.75
$ws.Cells.Item(27, 13).Value2 = -512.8
$ws.Cells.Item(27, 14).Value2 = -2163.75
# Row 68
$ws.Cells.Item(68, 8).Value2 = 1526.0834
$ws.Cells.Item(68, 9).Value2 = 1284
$ws.Cells.Item(68, 10).Value2 = 1768.1666
$ws.Cells.Item(68, 11).Value2 = 1284
$ws.Cells.Item(68, 12).Value2 = 1768.1666
$ws.Cells.Item(68, 13).Value2 = -535
$ws.Cells.Item(68, 14).Value2 = -3266.1666
# Row 71
$ws.Cells.Item(71, 8).Value2 = 1526.0834
$ws.Cells.Item(71, 9).Value2 = 1284
$ws.Cells.Item(71, 10).Value2 = 1768.1666
$ws.Cells.Item(71, 11).Value2 = 6420
$ws.Cells.Item(71, 12).Value2 = 8840.833000000001
$ws.Cells.Item(71, 13).Value2 = -2676
$ws.Cells.Item(71, 14).Value2 = -16328.833
# Row 82
$ws.Cells.Item(82, 8).Value2 = 2990
$ws.Cells.Item(82, 9).Value2 = 2958.5715
$ws.Cells.Item(82, 11).Value2 = 2958.5715
$ws.Cells.Item(82, 13).Value2 = -2597.5715
# Row 85
$ws.Cells.Item(85, 8).Value2 = 2990
$ws.Cells.Item(85, 9).Value2 = 2958.5715
$ws.Cells.Item(85, 11).Value2 = 2958.5715
$ws.Cells.Item(85, 13).Value2 = -1710.5715

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value2 = 833
$ws.Cells.Item(107, 9).Value2 = 749.5
$ws.Cells.Item(107, 10).Value2 = 1000
$ws.Cells.Item(107, 11).Value2 = 2248.5
$ws.Cells.Item(107, 12).Value2 = 3000
$ws.Cells.Item(107, 13).Value2 = -328.5
$ws.Cells.Item(107, 14).Value2 = -6840
